$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 268, shifting existing rows 268-308 down to 269-309.
$ws.Rows(268).Insert()

# Populate the newly inserted row 268 with the new record's data.
$ws.Range("A268").Value = 10
$ws.Range("B268").Value = "Vega Modelo de Temuco"
$ws.Range("C268").Value = "La Araucanía"
$ws.Range("D268").Value = 45180
$ws.Range("E268").Value = 9
$ws.Range("F268").Value = 100112012
$ws.Range("G268").Value = "Espinaca"
$ws.Range("H268").Value = "Sin especificar"
$ws.Range("I268").Value = "Primera"
$ws.Range("J268").Value = 40
$ws.Range("K268").Value = 12000
$ws.Range("L268").Value = 12000
$ws.Range("M268").Value = 12000
$ws.Range("N268").Value = "`$/docena de paquetes"
$ws.Range("O268").Value = "Región de La Araucanía"
$ws.Range("P268").Value = 1000
$ws.Range("Q268").Value = 12
$ws.Range("R268").Value = "Hortaliza"
